# ---------------------------------------------------------------------------
# 1) Slide 16: the single table on the slide switches from the deck's custom
#    "Table_0" style ({E7EA071D-BF84-4E9D-90A5-89AD963C2A59}) to the built-in
#    "Light Style 1" table style ({BE844B07-8FAF-4988-AC16-B0A37CD3118E}).
#    Table styles can't be set via a plain property assignment, so
#    Table.ApplyStyle(guid) is used (Table.StyleId is read-only).
# ---------------------------------------------------------------------------
$p = $ppt.ActivePresentation

$targetGuid = "{E7EA071D-BF84-4E9D-90A5-89AD963C2A59}"
$newGuid    = "{BE844B07-8FAF-4988-AC16-B0A37CD3118E}"

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTable) {
            $tbl = $shape.Table
            if ($tbl.StyleId -eq $targetGuid) {
                $tbl.ApplyStyle($newGuid)
            }
        }
    }
}

# ---------------------------------------------------------------------------
# 2) Theme swap: the presentation's active theme (the one behind the slide
#    master / Design 1, persisted as ppt/theme/theme2.xml) switches its
#    12-colour scheme from the custom "Integral" palette to the stock
#    "Office Theme" palette (fonts/format scheme are identical between the
#    two themes already, only the colours differ).
# ---------------------------------------------------------------------------
$master = $p.SlideMaster
$theme  = $master.Theme
$colorScheme = $theme.ThemeColorScheme

# Office Theme palette, in clrScheme document order:
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
$officeThemeColorsHex = @(
    "000000",  # dk1
    "FFFFFF",  # lt1
    "44546A",  # dk2
    "E7E6E6",  # lt2
    "5B9BD5",  # accent1
    "ED7D31",  # accent2
    "A5A5A5",  # accent3
    "FFC000",  # accent4
    "4472C4",  # accent5
    "70AD47",  # accent6
    "0563C1",  # hlink
    "954F72"   # folHlink
)

for ($i = 1; $i -le 12; $i++) {
    $hex = $officeThemeColorsHex[$i - 1]
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    $colorScheme.Colors($i).RGB = $r + ($g * 256) + ($b * 65536)
}
